$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instruments")
$ws.Range("B2").Value = "nhanes:00063"
